# feat: improve documentation for import process of samples and chemicals (#403)
#
# The "sample" import template had two redundant columns (AE "location" and
# AF "flash point") that duplicated the data already present in AB/AC using
# shorter "shelf N" labels instead of the full "room x- shelf N" / temperature
# text. These duplicate columns are removed, the now-wider "location" column
# (AB) is resized, the saved cursor position is updated, and a leftover
# duplicate "Normal" cell style (an artifact of the LibreOffice export) is
# cleaned up.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Drop the redundant duplicate "location"/"flash point" columns (AE/AF) ---
# Row 1 headers: keep the cells (and their style) but clear the text.
$ws.Range("AE1").ClearContents()
$ws.Range("AF1").ClearContents()

# Data rows 2-7: remove the duplicated "shelf N" / temperature values.
$ws.Range("AE2:AF2").ClearContents()
$ws.Range("AE3:AF3").ClearContents()
$ws.Range("AE4").ClearContents()
$ws.Range("AE5").ClearContents()
$ws.Range("AE6").ClearContents()
$ws.Range("AE7").ClearContents()

# --- Resize column AB ("location") now that it no longer needs to fit next
# to the removed duplicate column, and drop its old auto ("best fit") flag.
$ws.Columns("AB").ColumnWidth = 22.25

# --- Restore the author's last working selection ---
$ws.Range("AC20").Select()

# --- Remove the redundant "Normal" cell style left over from the original
# LibreOffice export, collapsing the style table back down to just the
# single base "Standard"/"Normal" style. ---
$wb.Styles("Normal").Delete()
